$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("Event") values for rows 2-25 (row 1 header already present)
$eventValues = @(0,0,0,1,0,0,0,0,0,0,0,1,0,0,0,0,0,0,1,1,1,0,0,0)

for ($i = 0; $i -lt $eventValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $eventValues[$i]
}

# Move the selection/scroll position as recorded in the saved view state
$ws.Range("E26").Select()
